# The workbook has 5 sheets: busbars, lines, loads, gens, trafos.
# On the "trafos" sheet a new column "v_base_kV" is inserted before the
# existing "V_SCH_pu" column (i.e. as the new column E), shifting the
# remaining header columns (V_SCH_pu .. tap_max) one column to the right.
$wb = $excel.ActiveWorkbook
$trafos = $wb.Worksheets.Item("trafos")

# Shift the existing header values in E1:N1 right into F1:O1 (from the
# last column down to column E) so nothing gets clobbered while copying.
for ($col = 14; $col -ge 5; $col--) {
    $srcCell = $trafos.Cells.Item(1, $col)
    $dstCell = $trafos.Cells.Item(1, $col + 1)
    $dstCell.Value = $srcCell.Value2
}

# Write the new header into the freed-up column E.
$trafos.Cells.Item(1, 5).Value = "v_base_kV"

# Make "trafos" the active sheet/tab and select its header row, matching
# the saved UI state of the edited workbook.
$trafos.Activate()
$trafos.Range("A1:O1").Select()
